# Rebuild the properties report with the new dataset (Huntsville Adventist Apartments
# repeated + the new South Bay Retirement Residence block), replacing the previous
# AHEPA 23-II Apartments report content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing hyperlinks + cell content/formatting before rebuilding the report
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# Row 1: single-letter column headers A..K
$ws.Cells.Item(1,1).Value = "A"
$ws.Cells.Item(1,2).Value = "B"
$ws.Cells.Item(1,3).Value = "C"
$ws.Cells.Item(1,4).Value = "D"
$ws.Cells.Item(1,5).Value = "E"
$ws.Cells.Item(1,6).Value = "F"
$ws.Cells.Item(1,7).Value = "G"
$ws.Cells.Item(1,8).Value = "H"
$ws.Cells.Item(1,9).Value = "I"
$ws.Cells.Item(1,10).Value = "J"
$ws.Cells.Item(1,11).Value = "K"

# Re-apply the bold, centered, bordered header style to row 1 (A1:K1)
$headerRng = $ws.Range("A1:K1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# --- Block starting at row 2 (AL property) ---
$ws.Cells.Item(2,1).Value = "#"
$ws.Cells.Item(2,2).Value = "Owner Organization Name"
$ws.Cells.Item(2,3).Value = "Property Name"
$ws.Cells.Item(2,4).Value = "Project Category"
$ws.Cells.Item(2,5).Value = "Owner Company Type"
$ws.Cells.Item(2,6).Value = "Projects Units"
$ws.Cells.Item(2,7).Value = "(Address) Line 1"
$ws.Cells.Item(2,8).Value = "(Address) City"
$ws.Cells.Item(2,9).Value = "(Address) State"
$ws.Cells.Item(2,10).Value = "(Address) Postal Code"
$ws.Cells.Item(2,11).Value = "ProPublica Link"
$ws.Cells.Item(3,1).Value = 800000304
$ws.Cells.Item(3,2).Value = "HUNTSVILLE ADVENTIST APARTMENTS, INC."
$ws.Cells.Item(3,3).Value = "HUNTSVILLE ADVENT APTS"
$ws.Cells.Item(3,4).Value = "Insured-Subsidized"
$ws.Cells.Item(3,5).Value = "Non-Profit"
$ws.Cells.Item(3,6).Value = 76
$ws.Cells.Item(3,7).Value = "3409 NW VENONA AVE"
$ws.Cells.Item(3,8).Value = "HUNTSVILLE"
$ws.Cells.Item(3,9).Value = "AL"
$ws.Cells.Item(3,11).Value = "https://projects.propublica.org/nonprofits/organizations/630878986"
$ws.Cells.Item(4,1).Value = "Contact Name"
$ws.Cells.Item(4,2).Value = "Phone"
$ws.Cells.Item(4,3).Value = "Email"
$ws.Cells.Item(4,4).Value = "Adress"
$ws.Cells.Item(4,5).Value = "Notes"
$ws.Cells.Item(5,1).Value = "Albert Dudley"
$ws.Cells.Item(6,1).Value = "Benjamin Jones"
$ws.Cells.Item(7,1).Value = "D Snell"
$ws.Cells.Item(8,1).Value = "Andryna Kuzmicic"
$ws.Cells.Item(9,1).Value = "Paul Gunn"

# --- Block starting at row 10 (AL property) ---
$ws.Cells.Item(10,1).Value = "#"
$ws.Cells.Item(10,2).Value = "Owner Organization Name"
$ws.Cells.Item(10,3).Value = "Property Name"
$ws.Cells.Item(10,4).Value = "Project Category"
$ws.Cells.Item(10,5).Value = "Owner Company Type"
$ws.Cells.Item(10,6).Value = "Projects Units"
$ws.Cells.Item(10,7).Value = "(Address) Line 1"
$ws.Cells.Item(10,8).Value = "(Address) City"
$ws.Cells.Item(10,9).Value = "(Address) State"
$ws.Cells.Item(10,10).Value = "(Address) Postal Code"
$ws.Cells.Item(10,11).Value = "ProPublica Link"
$ws.Cells.Item(11,1).Value = 800000304
$ws.Cells.Item(11,2).Value = "HUNTSVILLE ADVENTIST APARTMENTS, INC."
$ws.Cells.Item(11,3).Value = "HUNTSVILLE ADVENT APTS"
$ws.Cells.Item(11,4).Value = "Insured-Subsidized"
$ws.Cells.Item(11,5).Value = "Non-Profit"
$ws.Cells.Item(11,6).Value = 76
$ws.Cells.Item(11,7).Value = "3409 NW VENONA AVE"
$ws.Cells.Item(11,8).Value = "HUNTSVILLE"
$ws.Cells.Item(11,9).Value = "AL"
$ws.Cells.Item(11,11).Value = "https://projects.propublica.org/nonprofits/organizations/630878986"
$ws.Cells.Item(12,1).Value = "Contact Name"
$ws.Cells.Item(12,2).Value = "Phone"
$ws.Cells.Item(12,3).Value = "Email"
$ws.Cells.Item(12,4).Value = "Adress"
$ws.Cells.Item(12,5).Value = "Notes"
$ws.Cells.Item(13,1).Value = "Albert Dudley"
$ws.Cells.Item(14,1).Value = "Benjamin Jones"
$ws.Cells.Item(15,1).Value = "D Snell"
$ws.Cells.Item(16,1).Value = "Andryna Kuzmicic"
$ws.Cells.Item(17,1).Value = "Paul Gunn"

# --- Block starting at row 18 (AL property) ---
$ws.Cells.Item(18,1).Value = "#"
$ws.Cells.Item(18,2).Value = "Owner Organization Name"
$ws.Cells.Item(18,3).Value = "Property Name"
$ws.Cells.Item(18,4).Value = "Project Category"
$ws.Cells.Item(18,5).Value = "Owner Company Type"
$ws.Cells.Item(18,6).Value = "Projects Units"
$ws.Cells.Item(18,7).Value = "(Address) Line 1"
$ws.Cells.Item(18,8).Value = "(Address) City"
$ws.Cells.Item(18,9).Value = "(Address) State"
$ws.Cells.Item(18,10).Value = "(Address) Postal Code"
$ws.Cells.Item(18,11).Value = "ProPublica Link"
$ws.Cells.Item(19,1).Value = 800000304
$ws.Cells.Item(19,2).Value = "HUNTSVILLE ADVENTIST APARTMENTS, INC."
$ws.Cells.Item(19,3).Value = "HUNTSVILLE ADVENT APTS"
$ws.Cells.Item(19,4).Value = "Insured-Subsidized"
$ws.Cells.Item(19,5).Value = "Non-Profit"
$ws.Cells.Item(19,6).Value = 76
$ws.Cells.Item(19,7).Value = "3409 NW VENONA AVE"
$ws.Cells.Item(19,8).Value = "HUNTSVILLE"
$ws.Cells.Item(19,9).Value = "AL"
$ws.Cells.Item(19,11).Value = "https://projects.propublica.org/nonprofits/organizations/630878986"
$ws.Cells.Item(20,1).Value = "Contact Name"
$ws.Cells.Item(20,2).Value = "Phone"
$ws.Cells.Item(20,3).Value = "Email"
$ws.Cells.Item(20,4).Value = "Adress"
$ws.Cells.Item(20,5).Value = "Notes"
$ws.Cells.Item(21,1).Value = "Albert Dudley"
$ws.Cells.Item(22,1).Value = "Benjamin Jones"
$ws.Cells.Item(23,1).Value = "D Snell"
$ws.Cells.Item(24,1).Value = "Andryna Kuzmicic"
$ws.Cells.Item(25,1).Value = "Paul Gunn"

# --- Block starting at row 26 (AL property) ---
$ws.Cells.Item(26,1).Value = "#"
$ws.Cells.Item(26,2).Value = "Owner Organization Name"
$ws.Cells.Item(26,3).Value = "Property Name"
$ws.Cells.Item(26,4).Value = "Project Category"
$ws.Cells.Item(26,5).Value = "Owner Company Type"
$ws.Cells.Item(26,6).Value = "Projects Units"
$ws.Cells.Item(26,7).Value = "(Address) Line 1"
$ws.Cells.Item(26,8).Value = "(Address) City"
$ws.Cells.Item(26,9).Value = "(Address) State"
$ws.Cells.Item(26,10).Value = "(Address) Postal Code"
$ws.Cells.Item(26,11).Value = "ProPublica Link"
$ws.Cells.Item(27,1).Value = 800000304
$ws.Cells.Item(27,2).Value = "HUNTSVILLE ADVENTIST APARTMENTS, INC."
$ws.Cells.Item(27,3).Value = "HUNTSVILLE ADVENT APTS"
$ws.Cells.Item(27,4).Value = "Insured-Subsidized"
$ws.Cells.Item(27,5).Value = "Non-Profit"
$ws.Cells.Item(27,6).Value = 76
$ws.Cells.Item(27,7).Value = "3409 NW VENONA AVE"
$ws.Cells.Item(27,8).Value = "HUNTSVILLE"
$ws.Cells.Item(27,9).Value = "AL"
$ws.Cells.Item(27,11).Value = "https://projects.propublica.org/nonprofits/organizations/630878986"
$ws.Cells.Item(28,1).Value = "Contact Name"
$ws.Cells.Item(28,2).Value = "Phone"
$ws.Cells.Item(28,3).Value = "Email"
$ws.Cells.Item(28,4).Value = "Adress"
$ws.Cells.Item(28,5).Value = "Notes"
$ws.Cells.Item(29,1).Value = "Albert Dudley"
$ws.Cells.Item(30,1).Value = "Benjamin Jones"
$ws.Cells.Item(31,1).Value = "D Snell"
$ws.Cells.Item(32,1).Value = "Andryna Kuzmicic"
$ws.Cells.Item(33,1).Value = "Paul Gunn"

# --- Block starting at row 34 (CA property) ---
$ws.Cells.Item(34,1).Value = "#"
$ws.Cells.Item(34,2).Value = "Owner Organization Name"
$ws.Cells.Item(34,3).Value = "Property Name"
$ws.Cells.Item(34,4).Value = "Project Category"
$ws.Cells.Item(34,5).Value = "Owner Company Type"
$ws.Cells.Item(34,6).Value = "Projects Units"
$ws.Cells.Item(34,7).Value = "(Address) Line 1"
$ws.Cells.Item(34,8).Value = "(Address) City"
$ws.Cells.Item(34,9).Value = "(Address) State"
$ws.Cells.Item(34,10).Value = "(Address) Postal Code"
$ws.Cells.Item(34,11).Value = "ProPublica Link"
$ws.Cells.Item(35,1).Value = 800000067
$ws.Cells.Item(35,2).Value = "SOUTH BAY RETIREMENT RESIDENCE"
$ws.Cells.Item(35,3).Value = "South Bay Retirement Residence"
$ws.Cells.Item(35,4).Value = "202/811"
$ws.Cells.Item(35,5).Value = "Non-Profit"
$ws.Cells.Item(35,6).Value = 75
$ws.Cells.Item(35,7).Value = "1001 W CRESSEY ST"
$ws.Cells.Item(35,8).Value = "COMPTON"
$ws.Cells.Item(35,9).Value = "CA"
$ws.Cells.Item(35,11).Value = "https://projects.propublica.org/nonprofits/organizations/954321266"
$ws.Cells.Item(36,1).Value = "Contact Name"
$ws.Cells.Item(36,2).Value = "Phone"
$ws.Cells.Item(36,3).Value = "Email"
$ws.Cells.Item(36,4).Value = "Adress"
$ws.Cells.Item(36,5).Value = "Notes"
$ws.Cells.Item(37,1).Value = "Marie Hollis"
$ws.Cells.Item(38,1).Value = "Carlos Coates"
$ws.Cells.Item(39,1).Value = "Cheri L Blair"
$ws.Cells.Item(40,1).Value = "Mitzi Johnson"
$ws.Cells.Item(41,1).Value = "Rev Don Koepke"
$ws.Cells.Item(42,1).Value = "Anita McCrimon"
$ws.Cells.Item(43,1).Value = "Robert Covington"
$ws.Cells.Item(44,1).Value = "Martin T Laurent"
$ws.Cells.Item(45,1).Value = "Jean Enock Berus"
$ws.Cells.Item(46,1).Value = "Robert J Chillison II"

# --- Block starting at row 47 (AL property) ---
$ws.Cells.Item(47,1).Value = "#"
$ws.Cells.Item(47,2).Value = "Owner Organization Name"
$ws.Cells.Item(47,3).Value = "Property Name"
$ws.Cells.Item(47,4).Value = "Project Category"
$ws.Cells.Item(47,5).Value = "Owner Company Type"
$ws.Cells.Item(47,6).Value = "Projects Units"
$ws.Cells.Item(47,7).Value = "(Address) Line 1"
$ws.Cells.Item(47,8).Value = "(Address) City"
$ws.Cells.Item(47,9).Value = "(Address) State"
$ws.Cells.Item(47,10).Value = "(Address) Postal Code"
$ws.Cells.Item(47,11).Value = "ProPublica Link"
$ws.Cells.Item(48,1).Value = 800000304
$ws.Cells.Item(48,2).Value = "HUNTSVILLE ADVENTIST APARTMENTS, INC."
$ws.Cells.Item(48,3).Value = "HUNTSVILLE ADVENT APTS"
$ws.Cells.Item(48,4).Value = "Insured-Subsidized"
$ws.Cells.Item(48,5).Value = "Non-Profit"
$ws.Cells.Item(48,6).Value = 76
$ws.Cells.Item(48,7).Value = "3409 NW VENONA AVE"
$ws.Cells.Item(48,8).Value = "HUNTSVILLE"
$ws.Cells.Item(48,9).Value = "AL"
$ws.Cells.Item(48,11).Value = "https://projects.propublica.org/nonprofits/organizations/630878986"
$ws.Cells.Item(49,1).Value = "Contact Name"
$ws.Cells.Item(49,2).Value = "Phone"
$ws.Cells.Item(49,3).Value = "Email"
$ws.Cells.Item(49,4).Value = "Adress"
$ws.Cells.Item(49,5).Value = "Notes"
$ws.Cells.Item(50,1).Value = "Albert Dudley"
$ws.Cells.Item(51,1).Value = "Benjamin Jones"
$ws.Cells.Item(52,1).Value = "D Snell"
$ws.Cells.Item(53,1).Value = "Andryna Kuzmicic"
$ws.Cells.Item(54,1).Value = "Paul Gunn"

# Add the ProPublica hyperlinks (applied last so the K-column cells keep the
# dedicated "Hyperlink" style, matching the original report formatting)
$linkCell = $ws.Cells.Item(3,11)
$ws.Hyperlinks.Add($linkCell, "https://projects.propublica.org/nonprofits/organizations/630878986")
$linkCell.Style = "Hyperlink"
$linkCell = $ws.Cells.Item(11,11)
$ws.Hyperlinks.Add($linkCell, "https://projects.propublica.org/nonprofits/organizations/630878986")
$linkCell.Style = "Hyperlink"
$linkCell = $ws.Cells.Item(19,11)
$ws.Hyperlinks.Add($linkCell, "https://projects.propublica.org/nonprofits/organizations/630878986")
$linkCell.Style = "Hyperlink"
$linkCell = $ws.Cells.Item(27,11)
$ws.Hyperlinks.Add($linkCell, "https://projects.propublica.org/nonprofits/organizations/630878986")
$linkCell.Style = "Hyperlink"
$linkCell = $ws.Cells.Item(35,11)
$ws.Hyperlinks.Add($linkCell, "https://projects.propublica.org/nonprofits/organizations/954321266")
$linkCell.Style = "Hyperlink"
$linkCell = $ws.Cells.Item(48,11)
$ws.Hyperlinks.Add($linkCell, "https://projects.propublica.org/nonprofits/organizations/630878986")
$linkCell.Style = "Hyperlink"

